$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 26 (RM 232) - shifts rows up
$ws.Rows.Item(26).Delete()

# Now the old row 28 (SC 92) has become row 27; delete it too
$ws.Rows.Item(27).Delete()

# After the two row deletions, apply the remaining cell-level value edits.
# These correspond (by content) to: SC 5 (row26), SC101(row27), SC119(row29)
$ws.Range("D26").Value = ""
$ws.Range("D27").Value = -14.6
$ws.Range("D29").Value = ""
